# Atualizado por script em 05-11-2023 08:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 52-54: re-ordered match results (same matchday/date) ---
# Row 52 -> Decic vs Buducnost
$ws.Range("F52").Value = "Decic"
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = "Buducnost"
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 3.11
$ws.Range("K52").Value = "30/09/2023 12:43"
$ws.Range("L52").Value = 2.82
$ws.Range("M52").Value = "01/10/2023 17:39"
$ws.Range("N52").Value = 2.92
$ws.Range("O52").Value = "30/09/2023 12:43"
$ws.Range("P52").Value = 2.79
$ws.Range("Q52").Value = "01/10/2023 17:39"
$ws.Range("R52").Value = 2.22
$ws.Range("S52").Value = "30/09/2023 12:43"
$ws.Range("T52").Value = 2.79
$ws.Range("U52").Value = "01/10/2023 17:39"
$ws.Range("V52").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-buducnost/212AhUwR/"

# Row 53 -> Jezero vs Arsenal Tivat
$ws.Range("F53").Value = "Jezero"
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = "Arsenal Tivat"
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 2.08
$ws.Range("K53").Value = "30/09/2023 05:12"
$ws.Range("L53").Value = 2.41
$ws.Range("M53").Value = "01/10/2023 17:50"
$ws.Range("N53").Value = 2.84
$ws.Range("O53").Value = "30/09/2023 05:12"
$ws.Range("P53").Value = 2.84
$ws.Range("Q53").Value = "01/10/2023 17:50"
$ws.Range("R53").Value = 3.55
$ws.Range("S53").Value = "30/09/2023 05:12"
$ws.Range("T53").Value = 3.3
$ws.Range("U53").Value = "01/10/2023 17:50"
$ws.Range("V53").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-arsenal-tivat/fya2fj9E/"

# Row 54 -> Sutjeska vs Mornar Bar
$ws.Range("F54").Value = "Sutjeska"
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = "Mornar Bar"
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = 1.41
$ws.Range("K54").Value = "01/10/2023 12:43"
$ws.Range("L54").Value = 1.51
$ws.Range("M54").Value = "01/10/2023 17:46"
$ws.Range("N54").Value = 4.18
$ws.Range("O54").Value = "01/10/2023 12:43"
$ws.Range("P54").Value = 3.88
$ws.Range("Q54").Value = "01/10/2023 17:46"
$ws.Range("R54").Value = 7.18
$ws.Range("S54").Value = "01/10/2023 12:43"
$ws.Range("T54").Value = 6.55
$ws.Range("U54").Value = "01/10/2023 17:46"
$ws.Range("V54").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-mornar-bar/Gp0beWg8/"

# --- Rows 61 & 63: swapped match results (same matchday/date) ---
# Row 61 -> Arsenal Tivat vs Mornar Bar
$ws.Range("F61").Value = "Arsenal Tivat"
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = "Mornar Bar"
$ws.Range("I61").Value = 2
$ws.Range("J61").Value = 2.29
$ws.Range("K61").Value = "20/10/2023 02:12"
$ws.Range("L61").Value = 2.72
$ws.Range("M61").Value = "21/10/2023 14:43"
$ws.Range("N61").Value = 2.75
$ws.Range("O61").Value = "20/10/2023 02:12"
$ws.Range("P61").Value = 2.49
$ws.Range("Q61").Value = "21/10/2023 14:43"
$ws.Range("R61").Value = 3.2
$ws.Range("S61").Value = "20/10/2023 02:12"
$ws.Range("T61").Value = 3.32
$ws.Range("U61").Value = "21/10/2023 14:43"
$ws.Range("V61").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/arsenal-tivat-mornar-bar/rLJu89wE/"

# Row 63 -> Sutjeska vs Jedinstvo
$ws.Range("F63").Value = "Sutjeska"
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = "Jedinstvo"
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 1.29
$ws.Range("K63").Value = "20/10/2023 02:12"
$ws.Range("L63").Value = 1.42
$ws.Range("M63").Value = "21/10/2023 14:58"
$ws.Range("N63").Value = 4.78
$ws.Range("O63").Value = "20/10/2023 02:12"
$ws.Range("P63").Value = 4.4
$ws.Range("Q63").Value = "21/10/2023 14:58"
$ws.Range("R63").Value = 7.52
$ws.Range("S63").Value = "20/10/2023 02:12"
$ws.Range("T63").Value = 7.13
$ws.Range("U63").Value = "21/10/2023 14:58"
$ws.Range("V63").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-jedinstvo/MBIq7ThK/"

# --- Rows 68 & 69: swapped match results (same matchday/date) ---
# Row 68 -> Mladost DG vs Rudar
$ws.Range("F68").Value = "Mladost DG"
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = "Rudar"
$ws.Range("I68").Value = 1
$ws.Range("J68").Value = 2.17
$ws.Range("K68").Value = "27/10/2023 05:12"
$ws.Range("L68").Value = 1.53
$ws.Range("M68").Value = "28/10/2023 16:54"
$ws.Range("N68").Value = 2.98
$ws.Range("O68").Value = "27/10/2023 05:12"
$ws.Range("P68").Value = 3.93
$ws.Range("Q68").Value = "28/10/2023 16:55"
$ws.Range("R68").Value = 3.15
$ws.Range("S68").Value = "27/10/2023 05:12"
$ws.Range("T68").Value = 5.06
$ws.Range("U68").Value = "28/10/2023 16:55"
$ws.Range("V68").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-rudar/GS6Sb500/"

# Row 69 -> Mornar Bar vs Decic
$ws.Range("F69").Value = "Mornar Bar"
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = "Decic"
$ws.Range("I69").Value = 1
$ws.Range("J69").Value = 3.11
$ws.Range("K69").Value = "27/10/2023 05:12"
$ws.Range("L69").Value = 4.89
$ws.Range("M69").Value = "28/10/2023 16:07"
$ws.Range("N69").Value = 2.8
$ws.Range("O69").Value = "27/10/2023 05:12"
$ws.Range("P69").Value = 3.5
$ws.Range("Q69").Value = "28/10/2023 16:59"
$ws.Range("R69").Value = 2.3
$ws.Range("S69").Value = "27/10/2023 05:12"
$ws.Range("T69").Value = 1.65
$ws.Range("U69").Value = "28/10/2023 16:29"
$ws.Range("V69").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mornar-bar-decic/Y3HNaoof/"

# --- New row 71: Jezero vs Mladost DG ---
# Duplicate row 70's formatting (styles) down into the new row first.
$ws.Range("A70:V70").Copy($ws.Range("A71:V71"))

$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "montenegro"
$ws.Range("C71").Value = "prva-crnogorska-liga"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45234.66666666666
$ws.Range("F71").Value = "Jezero"
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = "Mladost DG"
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1.93
$ws.Range("K71").Value = "04/11/2023 02:13"
$ws.Range("L71").Value = 1.88
$ws.Range("M71").Value = "04/11/2023 15:45"
$ws.Range("N71").Value = 3.04
$ws.Range("O71").Value = "04/11/2023 02:13"
$ws.Range("P71").Value = 3.15
$ws.Range("Q71").Value = "04/11/2023 15:45"
$ws.Range("R71").Value = 3.74
$ws.Range("S71").Value = "04/11/2023 02:13"
$ws.Range("T71").Value = 4.49
$ws.Range("U71").Value = "04/11/2023 15:45"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-mladost-dg/lvEb1bMb/"
